$p = $ppt.ActivePresentation
$s = $p.Slides.Item(3)
$shape = $s.Shapes.Item(2)
$tr = $shape.TextFrame.TextRange
$para = $tr.Paragraphs(1, 1)

# The paragraph currently holds two runs:
#   Run 1: "Project-2"                                  (rPr without dirty="0")
#   Run 2: ": Health Nutrition ... inquiries:"           (rPr with dirty="0")
# The target text merges them into a single run (keeping run 2's formatting)
# with updated wording. Clear run 1's characters (removing that run), then
# rewrite the remaining run's text in place so formatting/run count match.
$firstRun = $para.Characters(1, 9)
$firstRun.Text = ""

$remaining = $para.Characters(1, $para.Length - 1)
$remaining.Text = "Project-2: Health Nutrition and population statistics for 14 governments at random, and applying machine learning to predict GDP expenditure on health by the government:"
